# Implement "My Appointments" page testing + "Upload files" functionality
# for the My Profile test case: a new "Update Image" / "ImageFile Name"
# pair of columns is inserted before the existing "Updated Gender" column,
# the My Appointments test case is marked executed, and the My Profile
# row gets sample upload data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before the old column V ("Updated Image URL"),
# shifting Updated Image URL/Updated Gender/Updated DOB one column to
# the right (V->W, W->X, X->Y).
$ws.Columns.Item(22).Insert()
$ws.Columns.Item(22).ColumnWidth = 20.83

# Row 1: stray value left over from editing the new columns.
$ws.Range("T1").Value = " "

# Row 2 (headers): label the two new columns.
$ws.Range("V2").Value = "Update Image"

# Row 5 sample data for the "My Profile Page Testing" case - set the
# uploaded file name before labeling the header (matches authoring order).
$ws.Range("W5").Value = "profile-pic.png"

$ws.Range("W2").Value = "ImageFile Name"

# Row 3: "My Appointments Page Testing" case now executed.
$ws.Range("C3").Value = "Yes"

# Row 5: mark the new "Update Image" step executed for the profile test.
$ws.Range("V5").Value = "Yes"

# The hidden _FilterDatabase name now also covers the header row.
$fd = $wb.Names.Item("Sheet1!_FilterDatabase")
$fd.RefersTo = "=Sheet1!`$C`$2:`$C`$5"

# Leave the selection where the author left it after the edit.
$ws.Range("H10").Select() | Out-Null
